$wb = $excel.ActiveWorkbook

# Delete column D (an empty/duplicate ID column) from the three "General" detail sheets.
# Excel will automatically shift remaining columns left, fix merged cells,
# hyperlinks, data validations and column widths.
$sheetNames = @("General - Documents", "General - Items", "General - Milestones")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("D:D").Delete()
}

# Make "General - Items" the active sheet/tab (mirrors the author ending their
# edit session there, after having started on "2 - Phase de lancement").
$itemsSheet = $wb.Worksheets.Item("General - Items")
$itemsSheet.Select()
$itemsSheet.Range("D2").Select()
